$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 476.58334
$ws.Range("I18").Value = 271.9
$ws.Range("J18").Value = 1500
$ws.Range("K18").Value = 271.9
$ws.Range("L18").Value = 1500
$ws.Range("M18").Value = 12.10000000000002
$ws.Range("N18").Value = -2068

$ws.Range("H129").Value = 838.5833
$ws.Range("I129").Value = 419.25
$ws.Range("J129").Value = 903.0961
$ws.Range("K129").Value = 1257.75
$ws.Range("L129").Value = 2709.2883
$ws.Range("M129").Value = 3742.25
$ws.Range("N129").Value = -12709.2883

$ws.Range("H137").Value = 3523.5
$ws.Range("I137").Value = 3393.7334
$ws.Range("J137").Value = 3700.4546
$ws.Range("K137").Value = 10181.2002
$ws.Range("L137").Value = 11101.3638
$ws.Range("M137").Value = -7631.200199999999
$ws.Range("N137").Value = -16201.3638

$ws.Range("H138").Value = 1881.9333
$ws.Range("I138").Value = 1630
$ws.Range("J138").Value = 2083.48
$ws.Range("K138").Value = 4890
$ws.Range("L138").Value = 6250.440000000001
$ws.Range("M138").Value = 250
$ws.Range("N138").Value = -16530.44

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1938.7
$ws.Range("I61").Value = 1422
$ws.Range("J61").Value = 4374.5713
$ws.Range("K61").Value = 1422
$ws.Range("L61").Value = 4374.5713
$ws.Range("M61").Value = -1210
$ws.Range("N61").Value = -4798.5713

$ws.Range("H74").Value = 1175.6957
$ws.Range("I74").Value = 1379.5625
$ws.Range("J74").Value = 709.7143
$ws.Range("K74").Value = 1379.5625
$ws.Range("L74").Value = 709.7143
$ws.Range("M74").Value = -505.5625
$ws.Range("N74").Value = -2457.7143

$ws.Range("H77").Value = 1175.6957
$ws.Range("I77").Value = 1379.5625
$ws.Range("J77").Value = 709.7143
$ws.Range("K77").Value = 6897.8125
$ws.Range("L77").Value = 3548.5715
$ws.Range("M77").Value = -2529.8125
$ws.Range("N77").Value = -12284.5715

$ws.Range("H125").Value = 69786
$ws.Range("J125").Value = 69786
$ws.Range("L125").Value = 69786
$ws.Range("N125").Value = -79626

$ws.Range("H132").Value = 2294.6604
$ws.Range("I132").Value = 1109.9
$ws.Range("J132").Value = 3840
$ws.Range("K132").Value = 3329.7
$ws.Range("L132").Value = 11520
$ws.Range("M132").Value = -799.7000000000003
$ws.Range("N132").Value = -16580

$ws.Range("H134").Value = 35174
$ws.Range("J134").Value = 35174
$ws.Range("L134").Value = 35174
$ws.Range("N134").Value = -45314

$ws.Range("H136").Value = 1938.7
$ws.Range("I136").Value = 1422
$ws.Range("J136").Value = 4374.5713
$ws.Range("K136").Value = 4266
$ws.Range("L136").Value = 13123.7139
$ws.Range("M136").Value = -1716
$ws.Range("N136").Value = -18223.7139

$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2501.5625
$ws.Range("I86").Value = 2751.5
$ws.Range("J86").Value = 1751.75
$ws.Range("K86").Value = 2751.5
$ws.Range("L86").Value = 1751.75
$ws.Range("M86").Value = -1628.5
$ws.Range("N86").Value = -3997.75

$ws.Range("H89").Value = 2501.5625
$ws.Range("I89").Value = 2751.5
$ws.Range("J89").Value = 1751.75
$ws.Range("K89").Value = 13757.5
$ws.Range("L89").Value = 8758.75
$ws.Range("M89").Value = -8141.5
$ws.Range("N89").Value = -19990.75

$ws.Range("H94").Value = 679.7826
$ws.Range("I94").Value = 674.3182
$ws.Range("J94").Value = 800
$ws.Range("K94").Value = 674.3182
$ws.Range("L94").Value = 800
$ws.Range("M94").Value = -223.3182
$ws.Range("N94").Value = -1702

$ws.Range("H107").Value = 1224.3334
$ws.Range("I107").Value = 842.28
$ws.Range("J107").Value = 6000
$ws.Range("K107").Value = 842.28
$ws.Range("L107").Value = 6000
$ws.Range("M107").Value = 1077.72
$ws.Range("N107").Value = -9840

$ws.Range("H112").Value = 32000
$ws.Range("J112").Value = 32000
$ws.Range("L112").Value = 32000
$ws.Range("N112").Value = -34954

$ws.Range("H134").Value = 2531.8386
$ws.Range("I134").Value = 1257.95
$ws.Range("K134").Value = 3773.85
$ws.Range("M134").Value = -1238.85

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1869.1578
$ws.Range("I31").Value = 1444.875
$ws.Range("J31").Value = 2867.4707
$ws.Range("K31").Value = 1444.875
$ws.Range("L31").Value = 2867.4707
$ws.Range("M31").Value = -1149.875
$ws.Range("N31").Value = -3457.4707

$ws.Range("H34").Value = 1869.1578
$ws.Range("I34").Value = 1444.875
$ws.Range("J34").Value = 2867.4707
$ws.Range("K34").Value = 1444.875
$ws.Range("L34").Value = 2867.4707
$ws.Range("M34").Value = -1242.875
$ws.Range("N34").Value = -3271.4707

$ws.Range("H58").Value = 1852.6818
$ws.Range("I58").Value = 1316.6111
$ws.Range("J58").Value = 2223.8076
$ws.Range("K58").Value = 1316.6111
$ws.Range("L58").Value = 2223.8076
$ws.Range("M58").Value = -1113.6111
$ws.Range("N58").Value = -2629.8076

$ws.Range("H132").Value = 2268.7058
$ws.Range("I132").Value = 1443.375
$ws.Range("J132").Value = 3002.3333
$ws.Range("K132").Value = 4330.125
$ws.Range("L132").Value = 9006.999899999999
$ws.Range("M132").Value = -1800.125
$ws.Range("N132").Value = -14066.9999

$ws.Range("H134").Value = 2310.111
$ws.Range("I134").Value = 2438.3333
$ws.Range("J134").Value = 2053.6667
$ws.Range("K134").Value = 7314.999899999999
$ws.Range("L134").Value = 6161.000100000001
$ws.Range("M134").Value = -4779.999899999999
$ws.Range("N134").Value = -11231.0001

$ws.Range("H136").Value = 1852.6818
$ws.Range("I136").Value = 1316.6111
$ws.Range("J136").Value = 2223.8076
$ws.Range("K136").Value = 3949.8333
$ws.Range("L136").Value = 6671.4228
$ws.Range("M136").Value = -1399.8333
$ws.Range("N136").Value = -11771.4228

$ws.Range("H140").Value = 43263.5
$ws.Range("J140").Value = 43263.5
$ws.Range("L140").Value = 43263.5
$ws.Range("N140").Value = -53623.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2259.55
$ws.Range("I132").Value = 1510.2142
$ws.Range("J132").Value = 4008
$ws.Range("K132").Value = 4530.642599999999
$ws.Range("L132").Value = 12024
$ws.Range("M132").Value = -2000.642599999999
$ws.Range("N132").Value = -17084

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2266.6667
$ws.Range("I61").Value = 1800
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 1800
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -1598
$ws.Range("N61").Value = -2904

$ws.Range("H113").Value = 2266.6667
$ws.Range("I113").Value = 1800
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 1800
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = 370
$ws.Range("N113").Value = -6840

$ws.Range("H132").Value = 7313.6523
$ws.Range("I132").Value = 9467.037
$ws.Range("J132").Value = 4253.579
$ws.Range("K132").Value = 28401.111
$ws.Range("L132").Value = 12760.737
$ws.Range("M132").Value = -25871.111
$ws.Range("N132").Value = -17820.737

$ws.Range("H136").Value = 9011640
$ws.Range("I136").Value = 2340.348
$ws.Range("J136").Value = 23812632
$ws.Range("K136").Value = 7021.044
$ws.Range("L136").Value = 71437896
$ws.Range("M136").Value = -4471.044
$ws.Range("N136").Value = -71442996

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 17485.715
$ws.Range("J94").Value = 17485.715
$ws.Range("L94").Value = 17485.715
$ws.Range("N94").Value = -19287.715

$ws.Range("H95").Value = 44000
$ws.Range("I95").Value = 1312
$ws.Range("J95").Value = 58229.332
$ws.Range("K95").Value = 1312
$ws.Range("L95").Value = 58229.332
$ws.Range("M95").Value = 1434
$ws.Range("N95").Value = -63721.332

$ws.Range("H98").Value = 24396.666
$ws.Range("J98").Value = 24396.666
$ws.Range("L98").Value = 24396.666
$ws.Range("N98").Value = -30386.666

$ws.Range("H102").Value = 31337
$ws.Range("J102").Value = 31337
$ws.Range("L102").Value = 31337
$ws.Range("N102").Value = -37827

$ws.Range("H103").Value = 16500
$ws.Range("J103").Value = 16500
$ws.Range("L103").Value = 16500
$ws.Range("N103").Value = -18844

$ws.Range("H110").Value = 36764.668
$ws.Range("J110").Value = 36764.668
$ws.Range("L110").Value = 36764.668
$ws.Range("N110").Value = -44944.668

$ws.Range("H114").Value = 37232
$ws.Range("J114").Value = 37232
$ws.Range("L114").Value = 37232
$ws.Range("N114").Value = -45910

$ws.Range("H132").Value = 1849.1714
$ws.Range("I132").Value = 1436.2222
$ws.Range("J132").Value = 3242.875
$ws.Range("K132").Value = 4308.6666
$ws.Range("L132").Value = 9728.625
$ws.Range("M132").Value = -1778.6666
$ws.Range("N132").Value = -14788.625

$ws.Range("H136").Value = 3054.3142
$ws.Range("I136").Value = 770.7619
$ws.Range("K136").Value = 2312.2857
$ws.Range("M136").Value = 237.7143000000001
